{"js": "// The document is a title paragraph followed by a table whose cells hold\n// short \"division fact\" strings. The edit replaces the title's date and\n// 25 of the table cell values. Since several cell values are duplicated\n// (e.g. \"91\u00f73=30, 1\" appears twice but maps to two different targets),\n// the replacement must be positional (document order), not a global\n// text find/replace.\n\nconst newValues = [\n  \"2025-11-30 Sunday\",\n  \"55\u00f76=9, 1\",\n  \"79\u00f79=8, 7\",\n  \"96\u00f79=10, 6\",\n  \"97\u00f79=10, 7\",\n  \"66\u00f75=13, 1\",\n  \"28\u00f79=3, 1\",\n  \"73\u00f78=9, 1\",\n  \"78\u00f76=13, 0\",\n  \"59\u00f75=11, 4\",\n  \"62\u00f79=6, 8\",\n  \"18\u00f73=6, 0\",\n  \"95\u00f73=31, 2\",\n  \"22\u00f74=5, 2\",\n  \"20\u00f77=2, 6\",\n  \"58\u00f73=19, 1\",\n  \"95\u00f79=10, 5\",\n  \"40\u00f77=5, 5\",\n  \"35\u00f77=5, 0\",\n  \"61\u00f78=7, 5\",\n  \"47\u00f74=11, 3\",\n  \"37\u00f78=4, 5\",\n  \"35\u00f72=17, 1\",\n  \"76\u00f79=8, 4\",\n  \"54\u00f72=27, 0\",\n  \"29\u00f77=4, 1\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet valueIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && valueIndex < newValues.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text !== \"\") {\n    para.insertText(newValues[valueIndex], \"Replace\");\n    valueIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a title paragraph followed by a 20-row x 5-column table.\n# Only 5 of the 20 rows (rows 1, 5, 9, 13, 17) actually hold \"division fact\"\n# text; the rest are blank spacer rows. The edit updates the title's date\n# and all 25 populated cells. Several source cell values are duplicated\n# (e.g. \"91\u00f73=30, 1\" appears in two different cells but maps to two\n# different targets), so the replacement is done positionally by\n# row/column (and by paragraph index for the title), not by text search.\n\n$d = $word.ActiveDocument\n\n# Title paragraph (the date line above the table).\n$d.Paragraphs.Item(1).Range.Text = \"2025-11-30 Sunday\"\n\n$tbl = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n  @(\"55\u00f76=9, 1\", \"79\u00f79=8, 7\", \"96\u00f79=10, 6\", \"97\u00f79=10, 7\", \"66\u00f75=13, 1\"),\n  @(\"28\u00f79=3, 1\", \"73\u00f78=9, 1\", \"78\u00f76=13, 0\", \"59\u00f75=11, 4\", \"62\u00f79=6, 8\"),\n  @(\"18\u00f73=6, 0\", \"95\u00f73=31, 2\", \"22\u00f74=5, 2\", \"20\u00f77=2, 6\", \"58\u00f73=19, 1\"),\n  @(\"95\u00f79=10, 5\", \"40\u00f77=5, 5\", \"35\u00f77=5, 0\", \"61\u00f78=7, 5\", \"47\u00f74=11, 3\"),\n  @(\"37\u00f78=4, 5\", \"35\u00f72=17, 1\", \"76\u00f79=8, 4\", \"54\u00f72=27, 0\", \"29\u00f77=4, 1\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n  $row = $dataRows[$i]\n  for ($col = 1; $col -le 5; $col++) {\n    $tbl.Cell($row, $col).Range.Text = $newValues[$i][$col - 1]\n  }\n}\n"}
